$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.310.81"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "3.083.50"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("E4").Value = "  +0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "233.31"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +7.23%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "624.54"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -5.13%  "
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("E9").Value = "  +0.09%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.728"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.62%  "
$ws.Range("E11").Value = "  -21.04%  "
$ws.Range("E12").Value = "  -3.37%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "36.36"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("E14").Value = "  -0.88%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.47"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.12%  "
$ws.Range("D16").Value = "90.201.40"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "3.652.52"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "3.088.79"
$ws.Range("E18").Value = "  -2.78%  "
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("E21").Value = "  -2.29%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "437.62"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.84%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.57"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +6.35%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "8.88"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "5.88"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.72%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.57"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.73%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "89.12"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "3.251.54"
$ws.Range("E29").Value = "  -2.66%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "9.52"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("E33").Value = "  -3.83%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.203"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +9.08%  "
$ws.Range("E35").Value = "  +1.21%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.154"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +6.26%  "
$ws.Range("E37").Value = "  +3.10%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "503.94"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("E39").Value = "  -0.68%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.99"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  -2.02%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0886"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.48%  "
$ws.Range("E43").Value = "  -3.65%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "22.19"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.51"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +52.03%  "
$ws.Range("E47").Value = "  -2.92%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "150.52"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.690"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.85%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("E51").Value = "  -0.86%  "